$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename player "Messi" -> "Montacer" (rows 2-6) and "Taha" -> "Yassine" (rows 7-13)
$ws.Range("B2:B6").Value2 = "Montacer"
$ws.Range("B7:B13").Value2 = "Yassine"

# Update the active selection on the sheet
$ws.Range("F12").Select()
